# Add a new row (96) of price data to the sheet, following the same
# pattern as the existing rows (date stored as literal text, the rest
# as numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 96

$ws.Cells.Item($row, 1).Value = "2024-11-06 00:00:00"
$ws.Cells.Item($row, 2).Value = 74950
$ws.Cells.Item($row, 3).Value = 10494.99
$ws.Cells.Item($row, 4).Value = 9287.6
$ws.Cells.Item($row, 5).Value = 7.1643
